$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'40.003.20"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +0.28%  "
$ws.Range("D3").Value = "'2.217.93"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -0.55%  "
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").Value = "'289.64"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -1.33%  "
$ws.Range("D6").Value = "'87.65"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +3.14%  "
$ws.Range("D7").Value = "'0.511"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -0.46%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("D9").Value = "'0.471"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +0.48%  "
$ws.Range("D10").Value = "'30.50"
$ws.Range("D10").ClearFormats()
$ws.Range("D11").Value = "'0.0777"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -2.75%  "
$ws.Range("E12").Value = "  +2.67%  "
$ws.Range("D13").Value = "'6.49"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +2.12%  "
$ws.Range("D14").Value = "'2.555.72"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -0.66%  "
$ws.Range("D15").Value = "'13.97"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -1.56%  "
$ws.Range("D16").Value = "'2.198.15"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -1.25%  "
$ws.Range("D17").Value = "'0.729"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +0.73%  "
$ws.Range("D18").Value = "'39.943.83"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +0.37%  "
$ws.Range("D19").Value = "'11.51"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +9.26%  "
$ws.Range("D20").Value = "'0.0$([char]0x2083)0882"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -0.91%  "
$ws.Range("E21").Value = "  +0.11%  "
$ws.Range("D22").Value = "'65.61"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +0.28%  "
$ws.Range("D23").Value = "'236.58"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +1.85%  "
$ws.Range("E25").Value = "  +1.28%  "
$ws.Range("E26").Value = "  -0.81%  "
$ws.Range("D27").Value = "'22.55"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -1.74%  "
$ws.Range("E28").Value = "  +0.01%  "
$ws.Range("D29").Value = "'9.21"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -0.14%  "
$ws.Range("D30").Value = "'156.02"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +0.75%  "
$ws.Range("D31").Value = "'31.81"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -3.30%  "
$ws.Range("E32").Value = "  -0.07%  "
$ws.Range("D33").Value = "'4.94"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +2.02%  "
$ws.Range("D34").Value = "'0.0718"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +1.56%  "
$ws.Range("E35").Value = "  +0.66%  "
$ws.Range("D36").Value = "'2.86"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +7.27%  "
$ws.Range("E37").Value = "  -0.18%  "
$ws.Range("D38").Value = "'15.83"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -3.67%  "
$ws.Range("D39").Value = "'0.0987"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +0.64%  "
$ws.Range("D40").Value = "'1.70"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +2.33%  "
$ws.Range("D41").Value = "'2.105.78"
$ws.Range("D41").ClearFormats()
$ws.Range("D42").Value = "'3.84"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +3.01%  "
$ws.Range("E43").Value = "  -1.42%  "
$ws.Range("E44").Value = "  +6.36%  "
$ws.Range("D45").Value = "'0.0267"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -0.98%  "
$ws.Range("D46").Value = "'17.35"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +7.17%  "
$ws.Range("E47").Value = "  +2.11%  "
$ws.Range("D48").Value = "'2.428.86"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -0.66%  "
$ws.Range("D49").Value = "'69.18"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -2.34%  "
$ws.Range("E50").Value = "  +0.14%  "
$ws.Range("D51").Value = "'88.61"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -0.28%  "
